$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.531.82"
$ws.Range("E2").Value = "  +1.28%  "
$ws.Range("D3").Value = "3.443.66"
$ws.Range("E3").Value = "  +0.72%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").Value = "'413.83"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.47%  "
$ws.Range("D6").Value = "'129.46"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.22%  "
$ws.Range("D7").Value = "'0.624"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.74%  "
$ws.Range("D8").Value = "'0.999"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Value = "'0.726"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.71%  "
$ws.Range("D10").Value = "'0.141"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.09%  "
$ws.Range("D11").Value = "'42.58"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.30%  "
$ws.Range("D12").Value = "'9.48"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.18%  "
$ws.Range("D13").Value = "'0.0000220"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +9.81%  "
$ws.Range("D14").Value = "3.986.56"
$ws.Range("E14").Value = "  +0.87%  "
$ws.Range("E15").Value = "  -0.21%  "
$ws.Range("D16").Value = "'20.49"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.57%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "3.454.93"
$ws.Range("E17").Value = "  +1.76%  "
$ws.Range("B18").Value = "Uniswap"
$ws.Range("C18").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D18").Value = "'12.91"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +6.03%  "
$ws.Range("E19").Value = "  +0.13%  "
$ws.Range("D20").Value = "62.515.64"
$ws.Range("E20").Value = "  +1.29%  "
$ws.Range("D21").Value = "'477.50"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +7.89%  "
$ws.Range("E22").Value = "  -0.17%  "
$ws.Range("E23").Value = "  +3.36%  "
$ws.Range("D24").Value = "'13.36"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.48%  "
$ws.Range("D25").Value = "'10.49"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +21.60%  "
$ws.Range("D26").Value = "'3.30"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.09%  "
$ws.Range("D27").Value = "'33.36"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.65%  "
$ws.Range("D28").Value = "'4.79"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.03%  "
$ws.Range("D29").Value = "'7.60"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.73%  "
$ws.Range("D30").Value = "'11.94"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.19%  "
$ws.Range("E31").Value = "  -2.36%  "
$ws.Range("D32").Value = "'0.166"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.25%  "
$ws.Range("E33").Value = "  -2.06%  "
$ws.Range("D34").Value = "'40.58"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.45%  "
$ws.Range("E35").Value = "  +0.16%  "
$ws.Range("D36").Value = "'58.65"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +10.36%  "
$ws.Range("D37").Value = "'0.0488"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.47%  "
$ws.Range("E38").Value = "  +0.06%  "
$ws.Range("D39").Value = "'3.01"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.89%  "
$ws.Range("D40").Value = "'0.324"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.16%  "
$ws.Range("B41").Value = "LidoDAOToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D41").Value = "'3.34"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.28%  "
$ws.Range("B42").Value = "Stellar"
$ws.Range("C42").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D42").Value = "'0.134"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.04%  "
$ws.Range("D43").Value = "'2.69"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +6.80%  "
$ws.Range("D44").Value = "'145.02"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.89%  "
$ws.Range("D45").Value = "'4.34"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.60%  "
$ws.Range("E46").Value = "  +4.56%  "
$ws.Range("D47").Value = "'2.41"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +12.44%  "
$ws.Range("D48").Value = "0.0₃0558"
$ws.Range("E48").Value = "  +38.49%  "
$ws.Range("D49").Value = "'16.34"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.15%  "
$ws.Range("D50").Value = "'22.41"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.01%  "
$ws.Range("D51").Value = "'110.98"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +8.07%  "
